{"js": "// Replace the worksheet date and every \"A\u00d7B=C\" answer cell with the\n// updated values from the commit. Each old value is unique in the\n// document, so a literal (non-wildcard) search-and-replace per pair is\n// unambiguous and keeps the original run formatting (font/size) intact.\nconst replacements = [\n  [\"2025-10-06 Monday\", \"2025-10-07 Tuesday\"],\n  [\"709\u00d76=4254\", \"283\u00d75=1415\"],\n  [\"556\u00d76=3336\", \"669\u00d72=1338\"],\n  [\"733\u00d75=3665\", \"433\u00d78=3464\"],\n  [\"258\u00d73=774\", \"617\u00d74=2468\"],\n  [\"464\u00d72=928\", \"576\u00d79=5184\"],\n  [\"559\u00d77=3913\", \"636\u00d73=1908\"],\n  [\"234\u00d73=702\", \"787\u00d75=3935\"],\n  [\"706\u00d77=4942\", \"898\u00d74=3592\"],\n  [\"749\u00d72=1498\", \"807\u00d74=3228\"],\n  [\"695\u00d75=3475\", \"511\u00d78=4088\"],\n  [\"668\u00d74=2672\", \"232\u00d72=464\"],\n  [\"745\u00d73=2235\", \"689\u00d78=5512\"],\n  [\"347\u00d72=694\", \"562\u00d75=2810\"],\n  [\"444\u00d77=3108\", \"523\u00d73=1569\"],\n  [\"474\u00d78=3792\", \"292\u00d77=2044\"],\n  [\"139\u00d76=834\", \"495\u00d74=1980\"],\n  [\"922\u00d78=7376\", \"698\u00d73=2094\"],\n  [\"816\u00d77=5712\", \"774\u00d76=4644\"],\n  [\"139\u00d79=1251\", \"128\u00d74=512\"],\n  [\"359\u00d74=1436\", \"950\u00d79=8550\"],\n  [\"593\u00d79=5337\", \"544\u00d76=3264\"],\n  [\"572\u00d79=5148\", \"293\u00d73=879\"],\n  [\"597\u00d79=5373\", \"389\u00d78=3112\"],\n  [\"712\u00d72=1424\", \"298\u00d79=2682\"],\n  [\"128\u00d78=1024\", \"656\u00d73=1968\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and every \"A\u00d7B=C\" answer cell with the\n# updated values from the commit. Each old value is unique in the\n# document, so Find/Replace (wdReplaceAll, one hit each) is unambiguous\n# and leaves the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-06 Monday\", \"2025-10-07 Tuesday\"),\n    @(\"709\u00d76=4254\", \"283\u00d75=1415\"),\n    @(\"556\u00d76=3336\", \"669\u00d72=1338\"),\n    @(\"733\u00d75=3665\", \"433\u00d78=3464\"),\n    @(\"258\u00d73=774\", \"617\u00d74=2468\"),\n    @(\"464\u00d72=928\", \"576\u00d79=5184\"),\n    @(\"559\u00d77=3913\", \"636\u00d73=1908\"),\n    @(\"234\u00d73=702\", \"787\u00d75=3935\"),\n    @(\"706\u00d77=4942\", \"898\u00d74=3592\"),\n    @(\"749\u00d72=1498\", \"807\u00d74=3228\"),\n    @(\"695\u00d75=3475\", \"511\u00d78=4088\"),\n    @(\"668\u00d74=2672\", \"232\u00d72=464\"),\n    @(\"745\u00d73=2235\", \"689\u00d78=5512\"),\n    @(\"347\u00d72=694\", \"562\u00d75=2810\"),\n    @(\"444\u00d77=3108\", \"523\u00d73=1569\"),\n    @(\"474\u00d78=3792\", \"292\u00d77=2044\"),\n    @(\"139\u00d76=834\", \"495\u00d74=1980\"),\n    @(\"922\u00d78=7376\", \"698\u00d73=2094\"),\n    @(\"816\u00d77=5712\", \"774\u00d76=4644\"),\n    @(\"139\u00d79=1251\", \"128\u00d74=512\"),\n    @(\"359\u00d74=1436\", \"950\u00d79=8550\"),\n    @(\"593\u00d79=5337\", \"544\u00d76=3264\"),\n    @(\"572\u00d79=5148\", \"293\u00d73=879\"),\n    @(\"597\u00d79=5373\", \"389\u00d78=3112\"),\n    @(\"712\u00d72=1424\", \"298\u00d79=2682\"),\n    @(\"128\u00d78=1024\", \"656\u00d73=1968\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        [ref]$oldText,   # FindText\n        [ref]$true,      # MatchCase\n        [ref]$false,     # MatchWholeWord\n        [ref]$false,     # MatchWildcards\n        [ref]$false,     # MatchSoundsLike\n        [ref]$false,     # MatchAllWordForms\n        [ref]$true,      # Forward\n        [ref]1,          # Wrap (wdFindContinue)\n        [ref]$false,     # Format\n        [ref]$newText,   # ReplaceWith\n        [ref]2           # Replace (wdReplaceAll)\n    )\n}\n"}
